{"js": "// Find the sentence ending \"... using the Plasmax medium.\" and update it to\n// \"... using the Plasmax, HPLM and RPMI media.\" (i.e. turn \"medium\" into\n// \", HPLM and RPMI media\").\nconst body = context.document.body;\n\nconst results = body.search(\" medium.\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the target text \" medium.\" in the document body.');\n}\n\n// There is exactly one occurrence in this document (right after \"Plasmax\").\nconst target = results.items[0];\ntarget.insertText(\", HPLM and RPMI media.\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Find the sentence ending \"... using the Plasmax medium.\" and update it to\n# \"... using the Plasmax, HPLM and RPMI media.\" (i.e. turn \"medium\" into\n# \", HPLM and RPMI media\").\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$found = $range.Find.Execute(\" medium.\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n\nif (-not $found) {\n    throw 'Could not find the target text \" medium.\" in the document.'\n}\n\n$range.Text = \", HPLM and RPMI media.\"\n"}
